# Change the table style applied to the (single) table on slides 14, 15
# and 16 from the custom "no style" table style to the built-in
# "No Style, Table Grid" style.
#
# In the PowerPoint UI this is done by selecting the table, opening the
# Table Design tab, and clicking the new style in the Table Styles
# gallery. Through the object model this is exposed as
# Table.ApplyStyle(styleId) (Table.Style is read-only and explicitly
# rejects direct assignment).

$p = $ppt.ActivePresentation

$newStyleId = "{D562C6D4-8E36-4443-8250-338C6E486F2A}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $null

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $tableShape = $shape
            break
        }
    }

    if ($tableShape -ne $null) {
        $tableShape.Table.ApplyStyle($newStyleId)
    }
}
